# Apply the cryptos list update (auto-generated from diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.259.66"
$ws.Range("E2").Value = "  +0.68%  "

$ws.Range("D3").Value = "1.863.83"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4679"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.59%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2857"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06538"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +14.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07905"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.18%  "

$ws.Range("D13").Value = "1.869.01"
$ws.Range("E13").Value = "  +0.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.164"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6817"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").Value = "30.260.86"
$ws.Range("E17").Value = "  +0.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.21%  "

$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.386"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.31%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007326"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.77%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.112.79"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.167"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.264"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.935"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.381"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09772"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.391"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.477"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.060"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04745"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.137"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7096"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.705"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01867"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.615"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.23%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.52%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.299"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.958"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8494"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4183"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9993"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "965.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.237"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.366"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05647"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.48%  "
